$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(15, 8).Value = 1238.0333  # H15
$ws.Cells.Item(15, 9).Value = 1238.0333  # I15
$ws.Cells.Item(15, 11).Value = 3714.0999  # K15
$ws.Cells.Item(15, 13).Value = -3545.0999  # M15

$ws.Cells.Item(137, 8).Value = 1402.5  # H137
$ws.Cells.Item(137, 9).Value = 1212.5  # I137
$ws.Cells.Item(137, 10).Value = 1627.0454  # J137
$ws.Cells.Item(137, 11).Value = 3637.5  # K137
$ws.Cells.Item(137, 12).Value = 4881.1362  # L137
$ws.Cells.Item(137, 13).Value = -1087.5  # M137
$ws.Cells.Item(137, 14).Value = -9981.136200000001  # N137

$ws.Cells.Item(138, 8).Value = 2754.51  # H138
$ws.Cells.Item(138, 9).Value = 1577.55  # I138
$ws.Cells.Item(138, 10).Value = 3048.75  # J138
$ws.Cells.Item(138, 11).Value = 4732.65  # K138
$ws.Cells.Item(138, 12).Value = 9146.25  # L138
$ws.Cells.Item(138, 13).Value = 407.3500000000004  # M138
$ws.Cells.Item(138, 14).Value = -19426.25  # N138

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 123900.65  # H32
$ws.Cells.Item(32, 9).Value = 122304.445  # I32
$ws.Cells.Item(32, 10).Value = 142257  # J32
$ws.Cells.Item(32, 11).Value = 122304.445  # K32
$ws.Cells.Item(32, 12).Value = 142257  # L32
$ws.Cells.Item(32, 13).Value = -122017.445  # M32
$ws.Cells.Item(32, 14).Value = -142831  # N32

$ws.Cells.Item(45, 8).Value = 1723.9688  # H45
$ws.Cells.Item(45, 9).Value = 1540.2727  # I45
$ws.Cells.Item(45, 10).Value = 2128.1  # J45
$ws.Cells.Item(45, 11).Value = 1540.2727  # K45
$ws.Cells.Item(45, 12).Value = 2128.1  # L45
$ws.Cells.Item(45, 13).Value = -1163.2727  # M45
$ws.Cells.Item(45, 14).Value = -2882.1  # N45

$ws.Cells.Item(61, 8).Value = 2023.1464  # H61
$ws.Cells.Item(61, 9).Value = 1654.0322  # I61
$ws.Cells.Item(61, 10).Value = 3167.4  # J61
$ws.Cells.Item(61, 11).Value = 1654.0322  # K61
$ws.Cells.Item(61, 12).Value = 3167.4  # L61
$ws.Cells.Item(61, 13).Value = -1442.0322  # M61
$ws.Cells.Item(61, 14).Value = -3591.4  # N61

$ws.Cells.Item(74, 8).Value = 40661.55  # H74
$ws.Cells.Item(74, 9).Value = 52720.027  # I74
$ws.Cells.Item(74, 10).Value = 1471.5  # J74
$ws.Cells.Item(74, 11).Value = 52720.027  # K74
$ws.Cells.Item(74, 12).Value = 1471.5  # L74
$ws.Cells.Item(74, 13).Value = -51846.027  # M74
$ws.Cells.Item(74, 14).Value = -3219.5  # N74

$ws.Cells.Item(77, 8).Value = 40661.55  # H77
$ws.Cells.Item(77, 9).Value = 52720.027  # I77
$ws.Cells.Item(77, 10).Value = 1471.5  # J77
$ws.Cells.Item(77, 11).Value = 263600.135  # K77
$ws.Cells.Item(77, 12).Value = 7357.5  # L77
$ws.Cells.Item(77, 13).Value = -259232.135  # M77
$ws.Cells.Item(77, 14).Value = -16093.5  # N77

$ws.Cells.Item(81, 8).Value = 250232000  # H81
$ws.Cells.Item(81, 9).Value = 60000  # I81
$ws.Cells.Item(81, 10).Value = 333622660  # J81
$ws.Cells.Item(81, 11).Value = 60000  # K81
$ws.Cells.Item(81, 12).Value = 333622660  # L81
$ws.Cells.Item(81, 13).Value = -59002  # M81
$ws.Cells.Item(81, 14).Value = -333624656  # N81

$ws.Cells.Item(84, 8).Value = 250232000  # H84
$ws.Cells.Item(84, 9).Value = 60000  # I84
$ws.Cells.Item(84, 10).Value = 333622660  # J84
$ws.Cells.Item(84, 11).Value = 180000  # K84
$ws.Cells.Item(84, 12).Value = 1000867980  # L84
$ws.Cells.Item(84, 13).Value = -175008  # M84
$ws.Cells.Item(84, 14).Value = -1000877964  # N84

$ws.Cells.Item(136, 8).Value = 2023.1464  # H136
$ws.Cells.Item(136, 9).Value = 1654.0322  # I136
$ws.Cells.Item(136, 10).Value = 3167.4  # J136
$ws.Cells.Item(136, 11).Value = 4962.096600000001  # K136
$ws.Cells.Item(136, 12).Value = 9502.200000000001  # L136
$ws.Cells.Item(136, 13).Value = -2412.096600000001  # M136
$ws.Cells.Item(136, 14).Value = -14602.2  # N136

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(134, 8).Value = 24727.91  # H134
$ws.Cells.Item(134, 9).Value = 1914.2424  # I134
$ws.Cells.Item(134, 10).Value = 87465.5  # J134
$ws.Cells.Item(134, 11).Value = 5742.7272  # K134
$ws.Cells.Item(134, 12).Value = 262396.5  # L134
$ws.Cells.Item(134, 13).Value = -3207.7272  # M134
$ws.Cells.Item(134, 14).Value = -267466.5  # N134

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 1069.13  # H31
$ws.Cells.Item(31, 9).Value = 975.1799999999999  # I31
$ws.Cells.Item(31, 10).Value = 1163.08  # J31
$ws.Cells.Item(31, 11).Value = 975.1799999999999  # K31
$ws.Cells.Item(31, 12).Value = 1163.08  # L31
$ws.Cells.Item(31, 13).Value = -680.1799999999999  # M31
$ws.Cells.Item(31, 14).Value = -1753.08  # N31

$ws.Cells.Item(34, 8).Value = 1069.13  # H34
$ws.Cells.Item(34, 9).Value = 975.1799999999999  # I34
$ws.Cells.Item(34, 10).Value = 1163.08  # J34
$ws.Cells.Item(34, 11).Value = 975.1799999999999  # K34
$ws.Cells.Item(34, 12).Value = 1163.08  # L34
$ws.Cells.Item(34, 13).Value = -773.1799999999999  # M34
$ws.Cells.Item(34, 14).Value = -1567.08  # N34

$ws.Cells.Item(58, 8).Value = 3330.16  # H58
$ws.Cells.Item(58, 9).Value = 836.6429000000001  # I58
$ws.Cells.Item(58, 10).Value = 6503.727  # J58
$ws.Cells.Item(58, 11).Value = 836.6429000000001  # K58
$ws.Cells.Item(58, 12).Value = 6503.727  # L58
$ws.Cells.Item(58, 13).Value = -633.6429000000001  # M58
$ws.Cells.Item(58, 14).Value = -6909.727  # N58

$ws.Cells.Item(107, 8).Value = 466.2143  # H107
$ws.Cells.Item(107, 9).Value = 365.1  # I107
$ws.Cells.Item(107, 11).Value = 365.1  # K107
$ws.Cells.Item(107, 13).Value = 1554.9  # M107

$ws.Cells.Item(132, 8).Value = 2383.081  # H132
$ws.Cells.Item(132, 9).Value = 1162.3158  # I132
$ws.Cells.Item(132, 10).Value = 3671.6667  # J132
$ws.Cells.Item(132, 11).Value = 3486.9474  # K132
$ws.Cells.Item(132, 12).Value = 11015.0001  # L132
$ws.Cells.Item(132, 13).Value = -956.9474  # M132
$ws.Cells.Item(132, 14).Value = -16075.0001  # N132

$ws.Cells.Item(134, 8).Value = 6757687.5  # H134
$ws.Cells.Item(134, 9).Value = 898.9516  # I134
$ws.Cells.Item(134, 10).Value = 41667760  # J134
$ws.Cells.Item(134, 11).Value = 2696.8548  # K134
$ws.Cells.Item(134, 12).Value = 125003280  # L134
$ws.Cells.Item(134, 13).Value = -161.8548000000001  # M134
$ws.Cells.Item(134, 14).Value = -125008350  # N134

$ws.Cells.Item(136, 8).Value = 3330.16  # H136
$ws.Cells.Item(136, 9).Value = 836.6429000000001  # I136
$ws.Cells.Item(136, 10).Value = 6503.727  # J136
$ws.Cells.Item(136, 11).Value = 2509.9287  # K136
$ws.Cells.Item(136, 12).Value = 19511.181  # L136
$ws.Cells.Item(136, 13).Value = 40.07129999999961  # M136
$ws.Cells.Item(136, 14).Value = -24611.181  # N136

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 6252.8423  # H5
$ws.Cells.Item(5, 9).Value = 621.4286  # I5
$ws.Cells.Item(5, 10).Value = 9537.833000000001  # J5
$ws.Cells.Item(5, 11).Value = 1864.2858  # K5
$ws.Cells.Item(5, 12).Value = 28613.499  # L5
$ws.Cells.Item(5, 13).Value = -1752.2858  # M5
$ws.Cells.Item(5, 14).Value = -28837.499  # N5

$ws.Cells.Item(68, 8).Value = 3150.889  # H68
$ws.Cells.Item(68, 9).Value = 15057.714  # I68
$ws.Cells.Item(68, 10).Value = 1662.5358  # J68
$ws.Cells.Item(68, 11).Value = 45173.142  # K68
$ws.Cells.Item(68, 12).Value = 4987.607400000001  # L68
$ws.Cells.Item(68, 13).Value = -44362.142  # M68
$ws.Cells.Item(68, 14).Value = -6609.607400000001  # N68

$ws.Cells.Item(71, 8).Value = 3150.889  # H71
$ws.Cells.Item(71, 9).Value = 15057.714  # I71
$ws.Cells.Item(71, 10).Value = 1662.5358  # J71
$ws.Cells.Item(71, 11).Value = 135519.426  # K71
$ws.Cells.Item(71, 12).Value = 14962.8222  # L71
$ws.Cells.Item(71, 13).Value = -131463.426  # M71
$ws.Cells.Item(71, 14).Value = -23074.8222  # N71

$ws.Cells.Item(107, 8).Value = 1220.6  # H107
$ws.Cells.Item(107, 9).Value = 1220.6  # I107
$ws.Cells.Item(107, 10).Value = 0  # J107
$ws.Cells.Item(107, 11).Value = 3661.8  # K107
$ws.Cells.Item(107, 12).Value = 0  # L107
$ws.Cells.Item(107, 13).Value = -1741.8  # M107
$ws.Cells.Item(107, 14).ClearContents()  # N107

$ws.Cells.Item(135, 8).Value = 6252.8423  # H135
$ws.Cells.Item(135, 9).Value = 621.4286  # I135
$ws.Cells.Item(135, 10).Value = 9537.833000000001  # J135
$ws.Cells.Item(135, 11).Value = 5592.8574  # K135
$ws.Cells.Item(135, 12).Value = 85840.497  # L135
$ws.Cells.Item(135, 13).Value = -3057.8574  # M135
$ws.Cells.Item(135, 14).Value = -90910.497  # N135

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(122, 8).Value = 2571.6052  # H122
$ws.Cells.Item(122, 9).Value = 2620.9565  # I122
$ws.Cells.Item(122, 10).Value = 2495.9333  # J122
$ws.Cells.Item(122, 11).Value = 7862.869499999999  # K122
$ws.Cells.Item(122, 12).Value = 7487.7999  # L122
$ws.Cells.Item(122, 13).Value = -5412.869499999999  # M122
$ws.Cells.Item(122, 14).Value = -12387.7999  # N122

$ws.Cells.Item(132, 8).Value = 29062.422  # H132
$ws.Cells.Item(132, 9).Value = 1752.591  # I132
$ws.Cells.Item(132, 10).Value = 66613.44  # J132
$ws.Cells.Item(132, 11).Value = 5257.772999999999  # K132
$ws.Cells.Item(132, 12).Value = 199840.32  # L132
$ws.Cells.Item(132, 13).Value = -2727.772999999999  # M132
$ws.Cells.Item(132, 14).Value = -204900.32  # N132

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(82, 8).Value = 1442.3214  # H82
$ws.Cells.Item(82, 9).Value = 1054.7142  # I82
$ws.Cells.Item(82, 10).Value = 1571.5238  # J82
$ws.Cells.Item(82, 11).Value = 1054.7142  # K82
$ws.Cells.Item(82, 12).Value = 1571.5238  # L82
$ws.Cells.Item(82, 13).Value = -693.7141999999999  # M82
$ws.Cells.Item(82, 14).Value = -2293.5238  # N82

$ws.Cells.Item(85, 8).Value = 1442.3214  # H85
$ws.Cells.Item(85, 9).Value = 1054.7142  # I85
$ws.Cells.Item(85, 10).Value = 1571.5238  # J85
$ws.Cells.Item(85, 11).Value = 1054.7142  # K85
$ws.Cells.Item(85, 12).Value = 1571.5238  # L85
$ws.Cells.Item(85, 13).Value = 193.2858000000001  # M85
$ws.Cells.Item(85, 14).Value = -4067.5238  # N85

$ws.Cells.Item(136, 8).Value = 436227.9  # H136
$ws.Cells.Item(136, 9).Value = 625881.0600000001  # I136
$ws.Cells.Item(136, 10).Value = 2735  # J136
$ws.Cells.Item(136, 11).Value = 1877643.18  # K136
$ws.Cells.Item(136, 12).Value = 8205  # L136
$ws.Cells.Item(136, 13).Value = -1875093.18  # M136
$ws.Cells.Item(136, 14).Value = -13305  # N136
